# Amend the "CAN" (H) column: two earlier polls' Canada-wide numbers were
# corrected, and a newly-added poll shifted the weighted/unweighted sample
# sizes. LIB..OTH (H4:H9) lose their explicit direct formatting (falls back
# to the workbook's default style) as part of the correction; nw/nu (H10:H11)
# keep their existing formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset direct number formatting on H4:H9 back to the sheet default, then
# write the corrected party vote-share figures for the CAN column.
$ws.Range("H4:H9").NumberFormat = "General"

$ws.Range("H4").Value = 35
$ws.Range("H5").Value = 29
$ws.Range("H6").Value = 23
$ws.Range("H7").Value = 7
$ws.Range("H8").Value = 5
$ws.Range("H9").Value = 2

# Updated weighted / unweighted sample sizes for the CAN column.
$ws.Range("H10").Value = 1223
$ws.Range("H11").Value = 1248

# Leave the cursor where the author left it after the last edit.
$ws.Range("H12").Select()
